$d = $word.ActiveDocument

# Locate the first occurrence of the run of text "Issue page" - this is the
# bullet right after "Create Issue (Button)" / "10 score" (Public Screens ->
# User Dashboard -> Issue page), and precedes "Edit Issue page".
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute("Issue page", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Re-defining the "_GoBack" bookmark moves it from wherever it currently
    # is (right after the "#/" run in the "Route: #/" paragraph) to wrap the
    # "Issue page" run instead - matching what real Word does when it
    # remembers the last edit location.
    $d.Bookmarks.Add("_GoBack", $rng)
}
